$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header changes ---
# F1 label: Revenue -> Investment
$ws.Range("F1").Value = "Investment"

# New columns O1, P1 - copy style from N1 (bold/border header style) then set text
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("N1").Copy()
$ws.Range("P1").PasteSpecial(-4122)

$ws.Range("O1").Value = "Land Investment"
$ws.Range("P1").Value = "Workforce Investment"

# --- Row 1 relabeling of investment-related headers (I1..N1 stay Water/Emission/Land Saving, add Workforce Saving at L1, shift Water/Emission Investment) ---
$ws.Range("L1").Value = "Workforce Saving"
$ws.Range("M1").Value = "Water Investment"
$ws.Range("N1").Value = "Emission Investment"

# --- Row 2 unit changes ---
$ws.Range("C2").Value = "M kSh/FU"
$ws.Range("E2").Value = "M kSh"
$ws.Range("F2").Value = "M kSh"
$ws.Range("G2").Value = "M kSh/FU"
$ws.Range("H2").Value = "years"
$ws.Range("J2").Value = "kton/FU"
$ws.Range("K2").Value = "M kSh/FU"
$ws.Range("L2").Value = "M kSh/FU"
$ws.Range("M2").Value = "m3/FU"
$ws.Range("N2").Value = "kton/FU"

# New columns O2, P2 - copy style from N2 then set text
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("N2").Copy()
$ws.Range("P2").PasteSpecial(-4122)

$ws.Range("O2").Value = "M kSh/FU"
$ws.Range("P2").Value = "M kSh/FU"

# --- Row 4 value changes ---
$ws.Range("F4").Value = 400.9907172536477
$ws.Range("G4").Value = -0.9994304059073329
$ws.Range("H4").Value = -401.219249367952
$ws.Range("I4").Value = -0.05114082011277787
$ws.Range("J4").Value = -0.01008209880092181
$ws.Range("K4").Value = -0.2510589539306238
$ws.Range("L4").Value = -0.1366259895730764
$ws.Range("M4").Value = 1.799733230887796
$ws.Range("N4").Value = 1.096701021939225
$ws.Range("O4").Value = 0.4893780273851007
$ws.Range("P4").Value = 20.56325181527063
